$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '53.919.96'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '2.248.14'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '492.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '126.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0948'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('E11').Value = '  +3.00%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '2.649.80'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').Value = '53.846.27'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '2.229.99'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.05%  '
$ws.Range('E19').Value = '  +3.25%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '298.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.55%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.86%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '2.357.26'
$ws.Range('E27').Value = '  +2.54%  '
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '165.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').Value = '0.0₃0673'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.884'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.23%  '
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('E40').Value = '  +2.84%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('E45').Value = '  +3.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.68'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '234.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('E51').Value = '  +0.26%  '
